# Insert a new weekly data row at the top of the "Mora" price series
# (row 28, the first data row for this market) and push the existing
# rows 28-72 down to 29-73, preserving all their values/formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(28).Insert()

$ws.Cells.Item(28, 1).Value  = 6
$ws.Cells.Item(28, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(28, 3).Value  = "Metropolitana"
$ws.Cells.Item(28, 4).Value  = 44581
$ws.Cells.Item(28, 5).Value  = 13
$ws.Cells.Item(28, 6).Value  = "Fruta"
$ws.Cells.Item(28, 7).Value  = 100101
$ws.Cells.Item(28, 8).Value  = "Berries"
$ws.Cells.Item(28, 9).Value  = 100101008
$ws.Cells.Item(28, 10).Value = "Mora"
$ws.Cells.Item(28, 11).Value = "Sin especificar"
$ws.Cells.Item(28, 12).Value = "Primera"
$ws.Cells.Item(28, 13).Value = 400
$ws.Cells.Item(28, 14).Value = 6000
$ws.Cells.Item(28, 15).Value = 6000
$ws.Cells.Item(28, 16).Value = 6000
$ws.Cells.Item(28, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(28, 18).Value = "Provincia de Linares"
$ws.Cells.Item(28, 19).Value = 3000
$ws.Cells.Item(28, 20).Value = 2
